$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("N2").Value = "2018-12-31 00:00:00"

$ws.Range("O2").Value = 51799452.62
$ws.Range("P2").Value = 699014994.91
$ws.Range("Q2").Value = 622031060.67
$ws.Range("R2").Value = 36.305857029
$ws.Range("S2").Value = 551870615.79
$ws.Range("T2").Value = 551870615.79
$ws.Range("U2").Value = 39.3236460557
$ws.Range("V2").Value = 6901953.28
$ws.Range("W2").Value = 58150126.03
$ws.Range("X2").Value = 176733.87
$ws.Range("Y2").Value = 75203664.34
$ws.Range("Z2").Value = 75096920.33
$ws.Range("AA2").Value = 19826084.06
$ws.Range("AG2").Value = 4931631.7
$ws.Range("AP2").Value = 37.7108986915
$ws.Range("AQ2").Value = 42.224296912202
$ws.Range("AR2").Value = 48.929613786699
$ws.Range("AS2").Value = 55753452.62
$ws.Range("AT2").Value = 39.522702054203
